$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceNone = 0 (we do the replace ourselves via Range.Text
# so the new text picks up the formatting of the *found* range instead of Word's
# "replacement inherits first character" quirk).

# 1) "NOTA PARA BG Nº XX/DP/2025." -> "NOTA PARA BG Nº {nota_bg}/DP/2025."
$rng = $d.Content
$rng.Find.Execute("XX", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "{nota_bg}"

# 2) "1º TEN QCOBM ALFREDO" -> "{posto_graduacao} {quadro} {nome_militar}"
$rng = $d.Content
$rng.Find.Execute("1º TEN QCOBM ALFREDO", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "{posto_graduacao} {quadro} {nome_militar}"

# 3) "4.10.2025" -> "{data_doacao}"
$rng = $d.Content
$rng.Find.Execute("4.10.2025", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "{data_doacao}"

# 4) "Enfª Antonia Liomar P. de Carvalho" -> "{atestador}"
$rng = $d.Content
$rng.Find.Execute("Enfª Antonia Liomar P. de Carvalho", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "{atestador}"

# 5) "479.807" -> "{numero_coren}"
$rng = $d.Content
$rng.Find.Execute("479.807", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "{numero_coren}"

# 6) "01.01.022101.033790/2025-23" -> "{numero_siged}"
$rng = $d.Content
$rng.Find.Execute("01.01.022101.033790/2025-23", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "{numero_siged}"

# 7) "8 de outubro de 2025" -> "{data_atual}"
$rng = $d.Content
$rng.Find.Execute("8 de outubro de 2025", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "{data_atual}"

# 8) drop the stray _GoBack bookmark left over from the last edit session
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
